# Refresh the crypto price/volume table (GitHub Actions scheduled update).
# Price cells in column D are stored as TEXT (values like "43.551.21" use
# dots as thousands separators and would otherwise be auto-coerced into a
# number by Excel's COM layer), so each D write is forced to Text via
# NumberFormat "@" before the value is set, then the style is reset back to
# "Normal" so no stray quote-prefix formatting is left on the cell - only
# the cell's value changes, matching the source data refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "43.551.21"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +2.11%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.245.76"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.30%  "
$ws.Range("E4").Value = "  -0.31%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "318.44"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.85%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "100.57"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +3.04%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.584"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +2.28%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.566"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.94%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "37.41"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.18%  "
$ws.Range("E11").Value = "  +0.27%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "7.75"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +2.32%  "
$ws.Range("E13").Value = "  +2.53%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.871"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.50%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "14.43"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +3.94%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.229.49"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.82%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "43.465.78"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +2.23%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "14.30"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +3.73%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.68"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.36%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0978"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +3.66%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "65.59"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.81%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "3.21"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.35%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "237.65"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.75%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.19"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +4.08%  "
$ws.Range("E25").Value = "  +0.53%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "4.06"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +2.89%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.11"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.08%  "
$ws.Range("E28").Value = "  +2.88%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "6.45"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.71%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "37.01"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +13.41%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "20.37"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.06%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.0876"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.55%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "158.93"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.87%  "
$ws.Range("E34").Value = "  +0.40%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.23"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +4.90%  "
$ws.Range("E36").Value = "  -0.68%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.91"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +6.01%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.42"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.56%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.105"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.10%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.74"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +7.03%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0324"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.46%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "14.55"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +22.11%  "
$ws.Range("E43").Value = "  -0.06%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.843.04"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +2.33%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.204"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.13%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "84.96"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -3.99%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "5.33"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.62%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "8.84"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +3.51%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "74.99"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.80%  "
$ws.Range("B50").Value = "MultiversX"
$ws.Range("C50").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "58.66"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.23%  "
$ws.Range("B51").Value = "Aave"
$ws.Range("C51").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "103.82"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.57%  "
